$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "department" column (C) previously held "BRANSON SCHOOL OF BUSINESS AND
# TECHNOLOGY" for every course row. Replace it with the more specific
# department names for each course.
$ws.Range("C2").Value = "Management"
$ws.Range("C3").Value = "Logistics"
$ws.Range("C4").Value = "Logistics"
